# Apply cryptocurrency price/volume updates to Sheet1 (cryptos.xlsx)
# Generated from the authoritative cell-level diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.138.31"
$ws.Range("E2").Value = "  -1.16%  "
$ws.Range("D3").Value = "1.862.39"
$ws.Range("E3").Value = "  -0.91%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9996"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7099"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.82%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "241.35"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.26%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3096"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.94%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07625"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.28%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.63"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.81%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08351"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.85%  "
$ws.Range("D12").Value = "1.865.02"
$ws.Range("E12").Value = "  -2.73%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.199"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.85%  "
$ws.Range("E14").Value = "  -3.35%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.08"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.36%  "
$ws.Range("D16").Value = "29.190.20"
$ws.Range("E16").Value = "  -1.08%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.901"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.85%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "242.59"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.07%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007800"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.93%  "
$ws.Range("D20").Value = "2.111.63"
$ws.Range("E20").Value = "  -2.44%  "
$ws.Range("E21").Value = "  -2.25%  "
$ws.Range("E22").Value = "  +0.05%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.872"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.65%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.9996"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.06%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1585"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.85%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "163.98"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.945"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.21%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.41"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.17%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.320"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.07%  "
$ws.Range("E30").Value = "  -0.05%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.387"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.16%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.247"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.97%  "
$ws.Range("E33").Value = "  -2.97%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7943"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +8.94%  "
$ws.Range("E35").Value = "  -2.37%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.161"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.39%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.684"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.36%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01844"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.47%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.696"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.46%  "
$ws.Range("D40").Value = "1.163.05"
$ws.Range("E40").Value = "  -5.24%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.227"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.26%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8911"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.35%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "72.78"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.04%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9995"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.06%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "102.66"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.17%  "
$ws.Range("D46").Value = "2.010.01"
$ws.Range("E46").Value = "  -1.72%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5185"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.57%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.775"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.29%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.273"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.90%  "
$ws.Range("B50").Value = "Frax"
$ws.Range("C50").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.003"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.44%  "
$ws.Range("B51").Value = "TheSandbox"
$ws.Range("C51").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4268"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.65%  "

Write-Host "Applied crypto list update: 63 plain + 36 text-forced cells"
